# Commit: "Support case insensitive lookup on Column C (regions)"
#
# Updates the sample data rows (A2:E3) with new entry/region words and
# extends the font formatting used for the remaining blank rows (A4:E8)
# from Calibri to Arial, which introduces a new font + cell style.
# Also restores the active cell selection left behind by the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - existing entry, words swapped out
$ws.Range("A2").Value = "escoba"
$ws.Range("B2").Value = "s. f. "
$ws.Range("C2").Value = "Sopeira"
$ws.Range("D2").Value = "escombra "
$ws.Range("E2").Value = "escoba"

# Row 3 - previously blank, now filled in with a second entry
$ws.Range("A3").Value = "escoba"
$ws.Range("B3").Value = "s. f. "
$ws.Range("C3").Value = "Sopeira, Tolba"
$ws.Range("D3").Value = "escombra"
$ws.Range("E3").Value = "balea"

# Remaining blank rows switch font from Calibri to Arial
$ws.Range("A4:E8").Font.Name = "Arial"

# Leave the same cell selected as in the authored workbook
$ws.Range("B56").Select()
